$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.112.19"
$ws.Range("E2").Value = "  -5.50%  "
$ws.Range("D3").Value = "3.701.69"
$ws.Range("E3").Value = "  -4.95%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'586.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").Value = "'182.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.48%  "
$ws.Range("D7").Value = "3.697.34"
$ws.Range("E7").Value = "  -4.84%  "
$ws.Range("D8").Value = "'0.627"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.48%  "
$ws.Range("D9").Value = "'0.998"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "'0.712"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.59%  "
$ws.Range("D11").Value = "'0.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.03%  "
$ws.Range("D12").Value = "'53.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "'0.0000291"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -10.50%  "
$ws.Range("D14").Value = "'10.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.98%  "
$ws.Range("D15").Value = "4.304.57"
$ws.Range("E15").Value = "  -4.60%  "
$ws.Range("D16").Value = "3.704.98"
$ws.Range("E16").Value = "  -4.78%  "
$ws.Range("D17").Value = "'19.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.38%  "
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").Value = "'12.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.34%  "
$ws.Range("D20").Value = "'1.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.08%  "
$ws.Range("D21").Value = "67.728.29"
$ws.Range("E21").Value = "  -5.67%  "
$ws.Range("D22").Value = "'407.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.98%  "
$ws.Range("D23").Value = "'4.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.93%  "
$ws.Range("D24").Value = "'88.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.33%  "
$ws.Range("D25").Value = "'3.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.13%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "'11.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'12.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.20%  "
$ws.Range("D28").Value = "'3.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.16%  "
$ws.Range("E29").Value = "  +1.84%  "
$ws.Range("D30").Value = "'9.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.12%  "
$ws.Range("D31").Value = "'32.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.59%  "
$ws.Range("D32").Value = "'7.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").Value = "'12.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.65%  "
$ws.Range("D34").Value = "'0.116"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.50%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'65.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.28%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "'606.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'43.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -13.35%  "
$ws.Range("D38").Value = "0.0₃0899"
$ws.Range("E38").Value = "  -9.01%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.399"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.72%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -5.20%  "
$ws.Range("D43").Value = "'2.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.05%  "
$ws.Range("D44").Value = "'3.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.46%  "
$ws.Range("D45").Value = "'3.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.51%  "
$ws.Range("D46").Value = "'0.0435"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.82%  "
$ws.Range("D47").Value = "'9.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.41%  "
$ws.Range("D48").Value = "2.822.31"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "'0.133"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.16%  "
$ws.Range("D50").Value = "'2.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.34%  "
$ws.Range("D51").Value = "'3.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.85%  "
